# "inserting data in hydration level" -- also reflects the incidental
# sheet-navigation/selection state from the same editing session
# (Home -> Sleep -> Hydration Level).
$wb = $excel.ActiveWorkbook

# The user had clicked on column C of the Home sheet at some point.
$home = $wb.Worksheets.Item("Home")
$home.Activate()
$home.Range("C1:C1048576").Select()

# Passed through the Sleep sheet (it was the previously active tab; it
# loses tabSelected once another sheet is activated).
$sleep = $wb.Worksheets.Item("Sleep")
$sleep.Activate()

# Landed on Hydration Level and entered the new readings.
$hydration = $wb.Worksheets.Item("Hydration Level")
$hydration.Activate()

$hydration.Range("A2").Value = 120
$hydration.Range("B2").Value = 120

$hydration.Range("A3").Value = 100
$hydration.Range("B3").Value = 120

$hydration.Range("A4").Value = 110
$hydration.Range("B4").Value = 120

$hydration.Range("A5").Value = 90
$hydration.Range("B5").Value = 120

$hydration.Range("A6").Value = 100
$hydration.Range("B6").Value = 120

$hydration.Range("A7").Value = 130
$hydration.Range("B7").Value = 120

$hydration.Range("A8").Value = 120
$hydration.Range("B8").Value = 120

$hydration.Range("A9").Value = 120
$hydration.Range("B9").Value = 120

$hydration.Range("A10").Value = 110
$hydration.Range("B10").Value = 120

$hydration.Range("A11").Value = 120
$hydration.Range("B11").Value = 120

# A couple of the entered cells ended up center-aligned.
$hydration.Range("A2").HorizontalAlignment = -4108
$hydration.Range("A7").HorizontalAlignment = -4108

# Final selection/cursor position left on B11.
$hydration.Range("B11").Select()
